$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 185 (old rows 185-212 shift down to 187-214)
$ws.Rows.Item(185).Insert()
$ws.Rows.Item(185).Insert()

# Fill new row 185 with the latest weekly data point
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44491
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 100112044
$ws.Range("G185").Value = "Perejil"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 65
$ws.Range("K185").Value = 4000
$ws.Range("L185").Value = 4000
$ws.Range("M185").Value = 4000
$ws.Range("N185").Value = "$/docena de atados (3 kilos)"
$ws.Range("O185").Value = "Provincia de Cautín"
$ws.Range("P185").Value = 1333
$ws.Range("Q185").Value = 3
$ws.Range("R185").Value = "Hortaliza"

# Fill new row 186 with the latest weekly data point
$ws.Range("A186").Value = 10
$ws.Range("B186").Value = "Vega Modelo de Temuco"
$ws.Range("C186").Value = "La Araucanía"
$ws.Range("D186").Value = 44491
$ws.Range("E186").Value = 9
$ws.Range("F186").Value = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 40
$ws.Range("K186").Value = 5000
$ws.Range("L186").Value = 5000
$ws.Range("M186").Value = 5000
$ws.Range("N186").Value = "$/docena de atados (3 kilos)"
$ws.Range("O186").Value = "Región del Maule"
$ws.Range("P186").Value = 1667
$ws.Range("Q186").Value = 3
$ws.Range("R186").Value = "Hortaliza"
